$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '66.876.48'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -0.34%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.521.85'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +1.33%  '
$ws.Range("E4").Value = '  +0.04%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '585.83'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.13%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '177.48'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '3.523.75'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +1.30%  '
$ws.Range("B9").Value = 'XRP'
$ws.Range("C9").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.600'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.83%  '
$ws.Range("E10").Value = '  -0.71%  '
$ws.Range("E11").Value = '  -1.88%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '4.131.58'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +1.39%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '30.67'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -3.43%  '
$ws.Range("E15").Value = '  -2.27%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '66.919.20'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.34%  '
$ws.Range("E17").Value = '  -0.79%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '3.523.63'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +1.52%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '6.12'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.73%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '14.08'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -1.03%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '381.76'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -1.44%  '
$ws.Range("E22").Value = '  -0.95%  '
$ws.Range("E23").Value = '  +0.26%  '
$ws.Range("E24").Value = '  +0.57%  '
$ws.Range("E25").Value = '  +0.64%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '71.71'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -2.88%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.0000122'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +1.44%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '9.94'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -3.72%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '0.174'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.14%  '
$ws.Range("E30").Value = '  -0.06%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '6.01'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -1.45%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '24.63'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +5.12%  '
$ws.Range("E33").Value = '  -1.50%  '
$ws.Range("E34").Value = '  -2.64%  '
$ws.Range("E35").Value = '  -0.03%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '7.26'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -0.88%  '
$ws.Range("E37").Value = '  -0.77%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '158.93'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -3.08%  '
$ws.Range("B39").Value = 'EnergySwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '29.20'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +12.19%  '
$ws.Range("B40").Value = 'Mantle'
$ws.Range("C40").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.892'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +2.84%  '
$ws.Range("E41").Value = '  -2.93%  '
$ws.Range("E42").Value = '  -3.42%  '
$ws.Range("E43").Value = '  -2.72%  '
$ws.Range("E44").Value = '  -1.78%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '2.733.61'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -3.10%  '
$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.0707'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -1.63%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '25.71'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -4.86%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '40.65'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -2.15%  '
$ws.Range("E49").Value = '  -0.28%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '328.55'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -1.34%  '
$ws.Range("E51").Value = '  -1.38%  '
